$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting the existing "age/level" column to E
$ws.Range("D1:D2").EntireColumn.Insert()

# New column D ("level") values
$ws.Range("D1").Value = 400
$ws.Range("D2").Value = 300

# Update selection to match target state
$ws.Range("H6").Select()
